$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176+ down to 177+
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new record's data
$ws.Range("A176").Value = 10
$ws.Range("B176").Value = "Vega Modelo de Temuco"
$ws.Range("C176").Value = "La Araucanía"
$ws.Range("D176").Value = 44460
$ws.Range("E176").Value = 9
$ws.Range("F176").Value = "Fruta"
$ws.Range("G176").Value = 100108
$ws.Range("H176").Value = "Tropicales y subtropicales"
$ws.Range("I176").Value = 100108005
$ws.Range("J176").Value = "Piña"
$ws.Range("K176").Value = "Caramelo"
$ws.Range("L176").Value = "Segunda"
$ws.Range("M176").Value = 30
$ws.Range("N176").Value = 22000
$ws.Range("O176").Value = 22000
$ws.Range("P176").Value = 22000
$ws.Range("Q176").Value = "$/caja 14 unidades"
$ws.Range("R176").Value = "Ecuador"
$ws.Range("S176").Value = 1571
$ws.Range("T176").Value = 14
